# repull data, push all data, mean calculation
# Update the "dSF" column (F) values to reflect the repulled/recalculated
# data for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    7  = -4
    9  = -2
    14 = -3
    15 = -1
    17 = -2
    22 = -2
    23 = -2
    25 = 4
    27 = 1
    29 = 2
    31 = 3
    33 = -1
    41 = 0
    42 = -1
    45 = 2
    49 = -1
    53 = 2
    56 = 1
    59 = 3
    65 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
